# Removed relevance in evaluation metric
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of data on the sheet and strip out any
# "Relevance: x/5" line from the Evaluation Score column (column C).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.Contains("Relevance")) {
        $lines = $val -split "`n"
        $keep = $lines | Where-Object { -not $_.StartsWith("Relevance") }
        $newVal = [string]::Join("`n", $keep)
        $cell.Value = $newVal
    }
}

# Update the sheet view scroll/selection state to match the saved file.
$ws.Range("D15").Select()
$av = $excel.ActiveWindow
$av.ScrollRow = 10
$av.ScrollColumn = 1
